$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $c = $ws.Range($cellAddr)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.704.09"
$ws.Range("E2").Value = "  +1.59%  "

Set-TextValue "D3" "1.880.57"
$ws.Range("E3").Value = "  +1.59%  "

Set-TextValue "D4" "1.007"
$ws.Range("E4").Value = "  +0.44%  "

Set-TextValue "D5" "332.93"
$ws.Range("E5").Value = "  +2.40%  "

Set-TextValue "D6" "1.006"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("E7").Value = "  +3.55%  "

Set-TextValue "D8" "0.3943"
$ws.Range("E8").Value = "  +1.62%  "

Set-TextValue "D9" "47.88"
$ws.Range("E9").Value = "  -0.79%  "

Set-TextValue "D10" "0.08072"
$ws.Range("E10").Value = "  +2.06%  "

Set-TextValue "D11" "1.027"
$ws.Range("E11").Value = "  +1.50%  "

Set-TextValue "D12" "22.19"
$ws.Range("E12").Value = "  +4.08%  "

Set-TextValue "D13" "1.884.76"
$ws.Range("E13").Value = "  +0.90%  "

Set-TextValue "D14" "5.984"
$ws.Range("E14").Value = "  +1.49%  "

Set-TextValue "D15" "7.128"
$ws.Range("E15").Value = "  -0.12%  "

Set-TextValue "D16" "1.009"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.00001048"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D18" "87.14"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D19" "0.06682"
$ws.Range("E19").Value = "  +1.30%  "

Set-TextValue "D20" "17.25"
$ws.Range("E20").Value = "  +0.79%  "

Set-TextValue "D21" "1.007"
$ws.Range("E21").Value = "  +0.42%  "

Set-TextValue "D22" "5.544"
$ws.Range("E22").Value = "  +1.09%  "

Set-TextValue "D23" "27.714.36"
$ws.Range("E23").Value = "  +1.58%  "

Set-TextValue "D24" "11.04"
$ws.Range("E24").Value = "  +1.96%  "

Set-TextValue "D25" "2.310"
$ws.Range("E25").Value = "  +0.72%  "

Set-TextValue "D26" "2.105.05"
$ws.Range("E26").Value = "  +1.10%  "

Set-TextValue "D27" "160.06"
$ws.Range("E27").Value = "  +3.82%  "

Set-TextValue "D28" "20.21"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("E29").Value = "  +2.45%  "

Set-TextValue "D30" "5.576"
$ws.Range("E30").Value = "  +2.41%  "

Set-TextValue "D31" "121.79"
$ws.Range("E31").Value = "  +0.56%  "

Set-TextValue "D32" "0.9834"
$ws.Range("E32").Value = "  +4.40%  "

Set-TextValue "D33" "0.09498"
$ws.Range("E33").Value = "  +1.89%  "

Set-TextValue "D34" "1.448"
$ws.Range("E34").Value = "  +0.37%  "

Set-TextValue "D35" "3.609"
$ws.Range("E35").Value = "  +0.63%  "

Set-TextValue "D36" "5.356"
$ws.Range("E36").Value = "  +2.20%  "

Set-TextValue "D37" "0.06142"
$ws.Range("E37").Value = "  +1.93%  "

Set-TextValue "D38" "0.02263"
$ws.Range("E38").Value = "  +1.76%  "

Set-TextValue "D39" "1.232"
$ws.Range("E39").Value = "  +2.15%  "

Set-TextValue "D40" "8.140"
$ws.Range("E40").Value = "  +0.95%  "

Set-TextValue "D41" "0.6000"
$ws.Range("E41").Value = "  +1.42%  "

Set-TextValue "D42" "0.1899"
$ws.Range("E42").Value = "  +0.99%  "

Set-TextValue "D43" "10.29"
$ws.Range("E43").Value = "  +1.58%  "

Set-TextValue "D44" "1.261"
$ws.Range("E44").Value = "  -1.51%  "

Set-TextValue "D45" "0.5723"
$ws.Range("E45").Value = "  +2.45%  "

$ws.Range("E46").Value = "  +0.74%  "

Set-TextValue "D47" "1.948"
$ws.Range("E47").Value = "  +2.15%  "

Set-TextValue "D48" "3.399"

Set-TextValue "D49" "0.06921"
$ws.Range("E49").Value = "  +2.88%  "

Set-TextValue "D50" "114.50"
$ws.Range("E50").Value = "  +6.24%  "

$ws.Range("E51").Value = "  +2.09%  "
